# Self perception Inventory.xlsx — fill in the Questionnaire answers
# (section totals of 10 points each). The "Grille d'évaluation" sheet
# pulls every one of these via formulas, so it (and the radar chart
# fed from it) recomputes automatically once these are written.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Questionnaire")

$answers = [ordered]@{
    "B4"  = 0
    "B5"  = 0
    "B6"  = 2
    "B8"  = 2
    "B9"  = 2
    "B10" = 2
    "B11" = 2
    "B15" = 1
    "B16" = 1
    "B17" = 1
    "B18" = 4
    "B19" = 1
    "B20" = 0
    "B21" = 2
    "B22" = 0
    "B26" = 2
    "B27" = 1
    "B28" = 2
    "B29" = 1
    "B30" = 1
    "B31" = 0
    "B32" = 2
    "B33" = 1
    "B37" = 0
    "B38" = 2
    "B39" = 3
    "B40" = 0
    "B41" = 0
    "B42" = 3
    "B43" = 0
    "B44" = 2
    "B48" = 3
    "B49" = 3
    "B50" = 0
    "B51" = 2
    "B52" = 0
    "B53" = 0
    "B54" = 0
    "B55" = 2
    "B59" = 2
    "B60" = 2
    "B61" = 1
    "B62" = 0
    "B63" = 2
    "B64" = 1
    "B65" = 2
    "B66" = 0
    "B70" = 3
    "B71" = 1
    "B72" = 1
    "B73" = 0
    "B74" = 3
    "B75" = 2
    "B76" = 0
    "B77" = 0
}

foreach ($cellRef in $answers.Keys) {
    $ws.Range($cellRef).Value = $answers[$cellRef]
}

# Recalculate everything so the "Grille d'évaluation" formulas pick up
# the new totals.
$excel.CalculateFull()

# Match the final view/selection state: Questionnaire scrolled down with
# B76 selected, and "Grille d'évaluation" as the active/selected sheet.
$ws.Range("B76").Select()

$wsGrille = $wb.Worksheets.Item("Grille d'évaluation")
$wsGrille.Activate()
$wsGrille.Range("P4:Q4").Select()
